# khurana_cv/data/training.xlsx update
#
# Fixes a typo in the "institution" column for the NAEP Data Training row
# (row 3): "American Institute for Research" -> "American Institutes for
# Research" (matches the real name of the organization, AIR).
# Also restores the on-disk selection/active-cell state recorded the last
# time the workbook was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = "US Department of Education and American Institutes for Research"

$ws.Range("G3").Select()
